$d = $word.ActiveDocument

# 1) "A coleta de dados referentes a este mapeamento foi feita" ->
#    "A coleta de dados referentes a este mapeamento em Jaraguá do Sul foi feita"
$d.Content.Find.Execute(
    "a este mapeamento foi feita",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "a este mapeamento em Jaraguá do Sul foi feita", 2)

# 2) "por órgãos públicos da Prefeitura de Jaraguá do Sul, que utilizou" ->
#    "por órgãos públicos do município, que utilizaram"
$d.Content.Find.Execute(
    "públicos da Prefeitura de Jaraguá do Sul, que utilizou",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "públicos do município, que utilizaram", 2)

# 3) "para a segurança da equipe. Em contrapartida" ->
#    "para a segurança da equipe (dados obtidos com um dos envolvidos no processo de
#     coleta de dados).  Em contrapartida"
$d.Content.Find.Execute(
    "para a segurança da equipe. Em contrapartida",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "para a segurança da equipe (dados obtidos com um dos envolvidos no processo de coleta de dados).  Em contrapartida", 2)
